# Add a new "x_brez_let" column (T) to the config table on every sheet,
# matching the commit "add x_sub_annual to config".
$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Range("T1").Value = "x_brez_let"
}

# Sheet32 (SURS--0300230S--B1GQ / SURS--0300230S--P3_S13 rows) also has its
# xmin/xmax date override moved from row 2 to row 3 (N column only).
$ws32 = $wb.Worksheets.Item(32)
$ws32.Range("N2").Copy()
$ws32.Range("N3").PasteSpecial(-4122)
$ws32.Range("M2").ClearContents()
$ws32.Range("N2").ClearContents()
$ws32.Range("N3").Value = 44562
$excel.CutCopyMode = $false

# Re-select the newly added header cell on every sheet and make Sheet1 the
# active tab (it takes over tabSelected from Sheet32).
for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Range("T1").Select()
}

$wb.Worksheets.Item(1).Activate()
$wb.Worksheets.Item(1).Range("T1").Select()
